# Update "想去人数" (F column) counts on several rows across sheets,
# as produced by the latest data refresh (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 575
$ws1.Range("F9").Value  = 2219
$ws1.Range("F10").Value = 746
$ws1.Range("F13").Value = 856
$ws1.Range("F18").Value = 27
$ws1.Range("F26").Value = 119
$ws1.Range("F32").Value = 153
$ws1.Range("F42").Value = 16

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value  = 472
$ws2.Range("F16").Value = 642
$ws2.Range("F19").Value = 472

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 2177

# Sheet "全部类型" (all types, aggregate of the above)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2177
$ws4.Range("F11").Value = 575
$ws4.Range("F15").Value = 746
$ws4.Range("F20").Value = 856
$ws4.Range("F24").Value = 472
$ws4.Range("F26").Value = 27
$ws4.Range("F33").Value = 119
$ws4.Range("F39").Value = 153
